$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 24, shifting existing rows 24-140 down to 25-141
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new data record
$ws.Cells.Item(24, 1).Value = 4
$ws.Cells.Item(24, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(24, 3).Value = "Los Lagos"

$dateCell = $ws.Cells.Item(24, 4)
$dateCell.Value = 45030
$dateCell.NumberFormat = $ws.Cells.Item(25, 4).NumberFormat

$ws.Cells.Item(24, 5).Value = 10
$ws.Cells.Item(24, 6).Value = 100112031
$ws.Cells.Item(24, 7).Value = "Poroto verde"
$ws.Cells.Item(24, 8).Value = "Magnum"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 40
$ws.Cells.Item(24, 11).Value = 26000
$ws.Cells.Item(24, 12).Value = 26000
$ws.Cells.Item(24, 13).Value = 26000
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 1040
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
